$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 16.48010064758079
$ws.Cells.Item(2, 3).Value = 7.950677712841619
$ws.Cells.Item(2, 4).Value = 8.186079208065285
$ws.Cells.Item(2, 5).Value = 12.59484330226041
$ws.Cells.Item(2, 6).Value = 34.55390627995787
$ws.Cells.Item(2, 9).Value = 25.51917288694298
$ws.Cells.Item(2, 10).Value = 9.809917668969986
$ws.Cells.Item(2, 12).Value = 11.31865242984793
$ws.Cells.Item(2, 13).Value = 16.4303851080165
$ws.Cells.Item(2, 14).Value = 18.95874407860833
$ws.Cells.Item(2, 15).Value = 26.52238024649509
$ws.Cells.Item(3, 2).Value = 16.10954675742117
$ws.Cells.Item(3, 3).Value = 7.674397866029614
$ws.Cells.Item(3, 4).Value = 8.191019248885196
$ws.Cells.Item(3, 5).Value = 12.62125984151247
$ws.Cells.Item(3, 6).Value = 34.62114515603822
$ws.Cells.Item(3, 9).Value = 25.61777460339745
$ws.Cells.Item(3, 10).Value = 9.824198551371303
$ws.Cells.Item(3, 12).Value = 11.32182694776391
$ws.Cells.Item(3, 13).Value = 16.3545126046496
$ws.Cells.Item(3, 14).Value = 19.01046269371206
$ws.Cells.Item(3, 15).Value = 26.58807750408736
$ws.Cells.Item(4, 2).Value = 15.87987837277332
$ws.Cells.Item(4, 3).Value = 7.49819936785744
$ws.Cells.Item(4, 4).Value = 8.194803789756705
$ws.Cells.Item(4, 5).Value = 12.63843509017408
$ws.Cells.Item(4, 6).Value = 34.67004136761454
$ws.Cells.Item(4, 9).Value = 25.68294926579025
$ws.Cells.Item(4, 10).Value = 9.833423671681864
$ws.Cells.Item(4, 12).Value = 11.32500197788311
$ws.Cells.Item(4, 13).Value = 16.30979423943689
$ws.Cells.Item(4, 14).Value = 19.04397701760033
$ws.Cells.Item(4, 15).Value = 26.63396365734895
$ws.Cells.Item(5, 2).Value = 15.78587755075016
$ws.Cells.Item(5, 3).Value = 7.424814990294608
$ws.Cells.Item(5, 4).Value = 8.196535475916594
$ws.Cells.Item(5, 5).Value = 12.64567496706303
$ws.Cells.Item(5, 6).Value = 34.69187786804912
$ws.Cells.Item(5, 9).Value = 25.71067273599438
$ws.Cells.Item(5, 10).Value = 9.837298127752005
$ws.Cells.Item(5, 12).Value = 11.32660501688427
$ws.Cells.Item(5, 13).Value = 16.29205349543292
$ws.Cells.Item(5, 14).Value = 19.05807767521973
$ws.Cells.Item(5, 15).Value = 26.65405526740384
$ws.Cells.Item(6, 2).Value = 15.77024799402681
$ws.Cells.Item(6, 3).Value = 7.412536150652738
$ws.Cells.Item(6, 4).Value = 8.196834478677168
$ws.Cells.Item(6, 5).Value = 12.64689170546495
$ws.Cells.Item(6, 6).Value = 34.69561911855067
$ws.Cells.Item(6, 9).Value = 25.71534649604456
$ws.Cells.Item(6, 10).Value = 9.837948442357739
$ws.Cells.Item(6, 12).Value = 11.32688990068148
$ws.Cells.Item(6, 13).Value = 16.28913715210651
$ws.Cells.Item(6, 14).Value = 19.06044588280567
$ws.Cells.Item(6, 15).Value = 26.65747550108946
$ws.Cells.Item(7, 2).Value = 15.87861212487007
$ws.Cells.Item(7, 3).Value = 7.497215990915305
$ws.Cells.Item(7, 4).Value = 8.194826376098083
$ws.Cells.Item(7, 5).Value = 12.63853175381774
$ws.Cells.Item(7, 6).Value = 34.67032812949095
$ws.Cells.Item(7, 9).Value = 25.68331844104708
$ws.Cells.Item(7, 10).Value = 9.833475457297689
$ws.Cells.Item(7, 12).Value = 11.32502234394193
$ws.Cells.Item(7, 13).Value = 16.309553012148
$ws.Cells.Item(7, 14).Value = 19.04416538775306
$ws.Cells.Item(7, 15).Value = 26.63422898446576
$ws.Cells.Item(8, 2).Value = 16.35285609306511
$ws.Cells.Item(8, 3).Value = 7.856816696285626
$ws.Cells.Item(8, 4).Value = 8.187626914685131
$ws.Cells.Item(8, 5).Value = 12.60375378617601
$ws.Cells.Item(8, 6).Value = 34.57550884316466
$ws.Cells.Item(8, 9).Value = 25.55220880344107
$ws.Cells.Item(8, 10).Value = 9.814747148099492
$ws.Cells.Item(8, 12).Value = 11.31949315226159
$ws.Cells.Item(8, 13).Value = 16.4038444226948
$ws.Cells.Item(8, 14).Value = 18.97621220899608
$ws.Cells.Item(8, 15).Value = 26.54387983664516
$ws.Cells.Item(9, 2).Value = 17.26034225288761
$ws.Cells.Item(9, 3).Value = 8.507319149653398
$ws.Cells.Item(9, 4).Value = 8.17944697107278
$ws.Cells.Item(9, 5).Value = 12.54310833109811
$ws.Cells.Item(9, 6).Value = 34.45007840577208
$ws.Cells.Item(9, 9).Value = 25.33189331479678
$ws.Cells.Item(9, 10).Value = 9.781628826093421
$ws.Cells.Item(9, 12).Value = 11.31833466758911
$ws.Cells.Item(9, 13).Value = 16.60298794424171
$ws.Cells.Item(9, 14).Value = 18.85686616881368
$ws.Cells.Item(9, 15).Value = 26.41082466532642
$ws.Cells.Item(10, 2).Value = 17.9064605399833
$ws.Cells.Item(10, 3).Value = 8.94889968884039
$ws.Cells.Item(10, 4).Value = 8.177025507401664
$ws.Cells.Item(10, 5).Value = 12.50312027911501
$ws.Cells.Item(10, 6).Value = 34.39494853986371
$ws.Cells.Item(10, 9).Value = 25.19249787620129
$ws.Cells.Item(10, 10).Value = 9.759474983064475
$ws.Cells.Item(10, 12).Value = 11.32332830699926
$ws.Cells.Item(10, 13).Value = 16.75717128931009
$ws.Cells.Item(10, 14).Value = 18.77759833899385
$ws.Cells.Item(10, 15).Value = 26.34009120552792
$ws.Cells.Item(11, 2).Value = 18.19459842408565
$ws.Cells.Item(11, 3).Value = 9.141381776478298
$ws.Cells.Item(11, 4).Value = 8.176695925054936
$ws.Cells.Item(11, 5).Value = 12.48591278723257
$ws.Cells.Item(11, 6).Value = 34.37792632886298
$ws.Cells.Item(11, 9).Value = 25.13397197706381
$ws.Cells.Item(11, 10).Value = 9.749865078093846
$ws.Cells.Item(11, 12).Value = 11.32685538155563
$ws.Cells.Item(11, 13).Value = 16.8288419850804
$ws.Cells.Item(11, 14).Value = 18.74335132351004
$ws.Cells.Item(11, 15).Value = 26.31380250129316
$ws.Cells.Item(12, 2).Value = 18.3027765107556
$ws.Cells.Item(12, 3).Value = 9.213026904177637
$ws.Cells.Item(12, 4).Value = 8.176681471940329
$ws.Cells.Item(12, 5).Value = 12.47953756306028
$ws.Cells.Item(12, 6).Value = 34.37263950004024
$ws.Cells.Item(12, 9).Value = 25.11251322630918
$ws.Cells.Item(12, 10).Value = 9.746293014973514
$ws.Cells.Item(12, 12).Value = 11.32837022271224
$ws.Cells.Item(12, 13).Value = 16.85618625334668
$ws.Cells.Item(12, 14).Value = 18.73064250295618
$ws.Cells.Item(12, 15).Value = 26.30469574863991
$ws.Cells.Item(13, 2).Value = 18.27952148022211
$ws.Cells.Item(13, 3).Value = 9.197652685109805
$ws.Cells.Item(13, 4).Value = 8.176679687390575
$ws.Cells.Item(13, 5).Value = 12.48090432412451
$ws.Cells.Item(13, 6).Value = 34.37372655845844
$ws.Cells.Item(13, 9).Value = 25.11710343460688
$ws.Cells.Item(13, 10).Value = 9.747059347566656
$ws.Cells.Item(13, 12).Value = 11.32803602564327
$ws.Cells.Item(13, 13).Value = 16.85028835742233
$ws.Cells.Item(13, 14).Value = 18.73336803400663
$ws.Cells.Item(13, 15).Value = 26.30661930590338
$ws.Cells.Item(14, 2).Value = 18.2035175644264
$ws.Cells.Item(14, 3).Value = 9.147301180600149
$ws.Cells.Item(14, 4).Value = 8.176692528189653
$ws.Cells.Item(14, 5).Value = 12.48538547319084
$ws.Cells.Item(14, 6).Value = 34.37746814671628
$ws.Cells.Item(14, 9).Value = 25.13219244205225
$ws.Cells.Item(14, 10).Value = 9.749569861236292
$ws.Cells.Item(14, 12).Value = 11.32697642450708
$ws.Cells.Item(14, 13).Value = 16.83108761695778
$ws.Cells.Item(14, 14).Value = 18.74230055947224
$ws.Cells.Item(14, 15).Value = 26.31303627643382
$ws.Cells.Item(15, 2).Value = 18.15683851540814
$ws.Cells.Item(15, 3).Value = 9.116296470202993
$ws.Cells.Item(15, 4).Value = 8.176714744893056
$ws.Cells.Item(15, 5).Value = 12.48814863790227
$ws.Cells.Item(15, 6).Value = 34.37991093535052
$ws.Cells.Item(15, 9).Value = 25.14152657700885
$ws.Cells.Item(15, 10).Value = 9.751116339870856
$ws.Cells.Item(15, 12).Value = 11.32635068699908
$ws.Cells.Item(15, 13).Value = 16.8193527084746
$ws.Cells.Item(15, 14).Value = 18.74780579613982
$ws.Cells.Item(15, 15).Value = 26.31707735709648
$ws.Cells.Item(16, 2).Value = 17.88750479056152
$ws.Cells.Item(16, 3).Value = 8.936148371525023
$ws.Cells.Item(16, 4).Value = 8.177062524780542
$ws.Cells.Item(16, 5).Value = 12.50426457109177
$ws.Cells.Item(16, 6).Value = 34.39622316446736
$ws.Cells.Item(16, 9).Value = 25.19642109529016
$ws.Cells.Item(16, 10).Value = 9.760112404650281
$ws.Cells.Item(16, 12).Value = 11.32312295612448
$ws.Cells.Item(16, 13).Value = 16.75251689257619
$ws.Cells.Item(16, 14).Value = 18.77987285923335
$ws.Cells.Item(16, 15).Value = 26.34192785465578
$ws.Cells.Item(17, 2).Value = 17.72072100722216
$ws.Cells.Item(17, 3).Value = 8.823455635401816
$ws.Cells.Item(17, 4).Value = 8.177473135347196
$ws.Cells.Item(17, 5).Value = 12.51440264656635
$ws.Cells.Item(17, 6).Value = 34.40829428570965
$ws.Cells.Item(17, 9).Value = 25.23134929699855
$ws.Cells.Item(17, 10).Value = 9.765750850163874
$ws.Cells.Item(17, 12).Value = 11.32146348186067
$ws.Cells.Item(17, 13).Value = 16.71189624650024
$ws.Cells.Item(17, 14).Value = 18.80000854652748
$ws.Cells.Item(17, 15).Value = 26.35868210609475
$ws.Cells.Item(18, 2).Value = 17.62425252193892
$ws.Cells.Item(18, 3).Value = 8.757849413730543
$ws.Cells.Item(18, 4).Value = 8.177782008301206
$ws.Cells.Item(18, 5).Value = 12.52032638018644
$ws.Cells.Item(18, 6).Value = 34.41599560700443
$ws.Cells.Item(18, 9).Value = 25.25189887726382
$ws.Cells.Item(18, 10).Value = 9.769038001849303
$ws.Cells.Item(18, 12).Value = 11.32062724749892
$ws.Cells.Item(18, 13).Value = 16.68867750408117
$ws.Cells.Item(18, 14).Value = 18.81176069349647
$ws.Cells.Item(18, 15).Value = 26.36887287222873
$ws.Cells.Item(19, 2).Value = 17.59150061982855
$ws.Cells.Item(19, 3).Value = 8.735502068244596
$ws.Cells.Item(19, 4).Value = 8.177899093269737
$ws.Cells.Item(19, 5).Value = 12.52234797118865
$ws.Cells.Item(19, 6).Value = 34.41873335903199
$ws.Cells.Item(19, 9).Value = 25.25893554825385
$ws.Cells.Item(19, 10).Value = 9.770158552564157
$ws.Cells.Item(19, 12).Value = 11.32036446168635
$ws.Cells.Item(19, 13).Value = 16.6808414588209
$ws.Cells.Item(19, 14).Value = 18.81576910033504
$ws.Cells.Item(19, 15).Value = 26.37241841909709
$ws.Cells.Item(20, 2).Value = 17.73853195620788
$ws.Cells.Item(20, 3).Value = 8.835533839476966
$ws.Cells.Item(20, 4).Value = 8.177421905570695
$ws.Cells.Item(20, 5).Value = 12.51331385310969
$ws.Cells.Item(20, 6).Value = 34.40693080369466
$ws.Cells.Item(20, 9).Value = 25.22758353404572
$ws.Cells.Item(20, 10).Value = 9.765146069408438
$ws.Cells.Item(20, 12).Value = 11.32162790647973
$ws.Cells.Item(20, 13).Value = 16.71620547424731
$ws.Cells.Item(20, 14).Value = 18.79784741450091
$ws.Cells.Item(20, 15).Value = 26.35684122193113
$ws.Cells.Item(21, 2).Value = 18.2258678643393
$ws.Cells.Item(21, 3).Value = 9.162124657035948
$ws.Cells.Item(21, 4).Value = 8.176685766931437
$ws.Cells.Item(21, 5).Value = 12.48406543172429
$ws.Cells.Item(21, 6).Value = 34.37633769188248
$ws.Cells.Item(21, 9).Value = 25.12774132621227
$ws.Cells.Item(21, 10).Value = 9.748830646527102
$ws.Cells.Item(21, 12).Value = 11.32728280204687
$ws.Cells.Item(21, 13).Value = 16.83672192998559
$ws.Cells.Item(21, 14).Value = 18.7396698163395
$ws.Cells.Item(21, 15).Value = 26.31112842610588
$ws.Cells.Item(22, 2).Value = 18.5388871375909
$ws.Cells.Item(22, 3).Value = 9.368308050046117
$ws.Cells.Item(22, 4).Value = 8.176847556117748
$ws.Cells.Item(22, 5).Value = 12.46577081682503
$ws.Cells.Item(22, 6).Value = 34.36309937030742
$ws.Cells.Item(22, 9).Value = 25.06659114075855
$ws.Cells.Item(22, 10).Value = 9.738557961955085
$ws.Cells.Item(22, 12).Value = 11.3320225437418
$ws.Cells.Item(22, 13).Value = 16.91666854816269
$ws.Cells.Item(22, 14).Value = 18.7031612243618
$ws.Cells.Item(22, 15).Value = 26.28619696778113
$ws.Cells.Item(23, 2).Value = 18.37235620694548
$ws.Cells.Item(23, 3).Value = 9.258939427578907
$ws.Cells.Item(23, 4).Value = 8.176702603290833
$ws.Cells.Item(23, 5).Value = 12.47546005130749
$ws.Cells.Item(23, 6).Value = 34.36954669571021
$ws.Cells.Item(23, 9).Value = 25.09885238035697
$ws.Cells.Item(23, 10).Value = 9.744005063996049
$ws.Cells.Item(23, 12).Value = 11.32939778642843
$ws.Cells.Item(23, 13).Value = 16.87389677606095
$ws.Cells.Item(23, 14).Value = 18.72250829994503
$ws.Cells.Item(23, 15).Value = 26.29905051954455
$ws.Cells.Item(24, 2).Value = 17.7304814319772
$ws.Cells.Item(24, 3).Value = 8.830075827489329
$ws.Cells.Item(24, 4).Value = 8.177444839727945
$ws.Cells.Item(24, 5).Value = 12.51380580000095
$ws.Cells.Item(24, 6).Value = 34.40754486202016
$ws.Cells.Item(24, 9).Value = 25.22928457498143
$ws.Cells.Item(24, 10).Value = 9.765419348929989
$ws.Cells.Item(24, 12).Value = 11.32155320301705
$ws.Cells.Item(24, 13).Value = 16.714256851858
$ws.Cells.Item(24, 14).Value = 18.79882391438778
$ws.Cells.Item(24, 15).Value = 26.35767174599768
$ws.Cells.Item(25, 2).Value = 17.01799005618252
$ws.Cells.Item(25, 3).Value = 8.337530140963821
$ws.Cells.Item(25, 4).Value = 8.181027323145754
$ws.Cells.Item(25, 5).Value = 12.55870966525312
$ws.Cells.Item(25, 6).Value = 34.47751659659405
$ws.Cells.Item(25, 9).Value = 25.38755157784107
$ws.Cells.Item(25, 10).Value = 9.790204171981399
$ws.Cells.Item(25, 12).Value = 11.31761713263736
$ws.Cells.Item(25, 13).Value = 16.54767158026237
$ws.Cells.Item(25, 14).Value = 18.88766998462999
$ws.Cells.Item(25, 15).Value = 26.4420823866872
